$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1023.7826
$ws.Range("I40").Value = 1005.375
$ws.Range("K40").Value = 1005.375
$ws.Range("M40").Value = -830.375

$ws.Range("H98").Value = 4533
$ws.Range("I98").Value = 4063
$ws.Range("K98").Value = 4063
$ws.Range("M98").Value = -2565

$ws.Range("H112").Value = 1989.9756
$ws.Range("J112").Value = 2041.8158
$ws.Range("L112").Value = 6125.4474
$ws.Range("N112").Value = -8341.447400000001

$ws.Range("H122").Value = 4533
$ws.Range("I122").Value = 4063
$ws.Range("K122").Value = 12189
$ws.Range("M122").Value = -9739

$ws.Range("H127").Value = 2087.9473
$ws.Range("I127").Value = 1762.5
$ws.Range("K127").Value = 5287.5
$ws.Range("M127").Value = -327.5

$ws.Range("H131").Value = 2122.1667
$ws.Range("I131").Value = 727.13336
$ws.Range("K131").Value = 2181.40008
$ws.Range("M131").Value = 2858.59992

$ws.Range("H135").Value = 38462076
$ws.Range("I135").Value = 520.125
$ws.Range("J135").Value = 500000740
$ws.Range("K135").Value = 4681.125
$ws.Range("L135").Value = 4500006660
$ws.Range("M135").Value = -2146.125
$ws.Range("N135").Value = -4500011730

$ws.Range("H138").Value = 1576.5476
$ws.Range("I138").Value = 1266.7213
$ws.Range("K138").Value = 3800.1639
$ws.Range("M138").Value = 1339.8361

$ws.Range("H139").Value = 50908.145
$ws.Range("J139").Value = 50908.145
$ws.Range("L139").Value = 50908.145
$ws.Range("N139").Value = -61188.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2778277.5
$ws.Range("I2").Value = 2778277.5
$ws.Range("K2").Value = 2778277.5
$ws.Range("M2").Value = -2778164.5

$ws.Range("H32").Value = 3759.9583
$ws.Range("I32").Value = 3063.0793
$ws.Range("K32").Value = 3063.0793
$ws.Range("M32").Value = -2776.0793

$ws.Range("H61").Value = 33334984
$ws.Range("I61").Value = 22728288
$ws.Range("K61").Value = 22728288
$ws.Range("M61").Value = -22728076

$ws.Range("H63").Value = 9250.833000000001
$ws.Range("I63").Value = 8101
$ws.Range("J63").Value = 15000
$ws.Range("K63").Value = 8101
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -7415

$ws.Range("H66").Value = 9250.833000000001
$ws.Range("I66").Value = 8101
$ws.Range("J66").Value = 15000
$ws.Range("K66").Value = 40505
$ws.Range("L66").Value = 75000
$ws.Range("M66").Value = -37073

$ws.Range("H74").Value = 1407.5122
$ws.Range("I74").Value = 1174.4138
$ws.Range("K74").Value = 1174.4138
$ws.Range("M74").Value = -300.4138

$ws.Range("H77").Value = 1407.5122
$ws.Range("I77").Value = 1174.4138
$ws.Range("K77").Value = 5872.069
$ws.Range("M77").Value = -1504.069

$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000

$ws.Range("H116").Value = 2778277.5
$ws.Range("I116").Value = 2778277.5
$ws.Range("K116").Value = 2778277.5
$ws.Range("M116").Value = -2775983.5

$ws.Range("H132").Value = 1414.317
$ws.Range("I132").Value = 1178.9032
$ws.Range("K132").Value = 3536.7096
$ws.Range("M132").Value = -1006.7096

$ws.Range("H136").Value = 33334984
$ws.Range("I136").Value = 22728288
$ws.Range("K136").Value = 68184864
$ws.Range("M136").Value = -68182314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2778277.5
$ws.Range("I3").Value = 2778277.5
$ws.Range("K3").Value = 2778277.5
$ws.Range("M3").Value = -2778163.5

$ws.Range("H107").Value = 1870.4
$ws.Range("I107").Value = 1863
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 1863
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 57

$ws.Range("H134").Value = 5183
$ws.Range("I134").Value = 5183
$ws.Range("K134").Value = 15549
$ws.Range("M134").Value = -13014

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3176.3157
$ws.Range("I122").Value = 2375.9167
$ws.Range("K122").Value = 7127.750100000001
$ws.Range("M122").Value = -4677.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1001.9286
$ws.Range("I122").Value = 777.8
$ws.Range("J122").Value = 1050.6522
$ws.Range("K122").Value = 7000.2
$ws.Range("L122").Value = 9455.8698
$ws.Range("M122").Value = -4550.2
$ws.Range("N122").Value = -14355.8698

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0

$ws.Range("H122").Value = 1115.6666
$ws.Range("J122").Value = 1492
$ws.Range("L122").Value = 4476
$ws.Range("N122").Value = -9376

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2987.8823
$ws.Range("I7").Value = 2634
$ws.Range("J7").Value = 3180.9092
$ws.Range("K7").Value = 2634
$ws.Range("L7").Value = 3180.9092
$ws.Range("M7").Value = -2522
$ws.Range("N7").Value = -3404.9092

$ws.Range("H61").Value = 2424.1875
$ws.Range("I61").Value = 2056.2144
$ws.Range("K61").Value = 2056.2144
$ws.Range("M61").Value = -1854.2144

$ws.Range("H113").Value = 2424.1875
$ws.Range("I113").Value = 2056.2144
$ws.Range("K113").Value = 2056.2144
$ws.Range("M113").Value = 113.7856000000002

$ws.Range("H126").Value = 2987.8823
$ws.Range("I126").Value = 2634
$ws.Range("J126").Value = 3180.9092
$ws.Range("K126").Value = 7902
$ws.Range("L126").Value = 9542.7276
$ws.Range("M126").Value = -5432
$ws.Range("N126").Value = -14482.7276

$ws.Range("H136").Value = 2080.7908
$ws.Range("I136").Value = 1364.5161
$ws.Range("J136").Value = 3931.1667
$ws.Range("K136").Value = 4093.5483
$ws.Range("L136").Value = 11793.5001
$ws.Range("M136").Value = -1543.5483
$ws.Range("N136").Value = -16893.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1212.9365
$ws.Range("I132").Value = 933.56604
$ws.Range("J132").Value = 2693.6
$ws.Range("K132").Value = 2800.69812
$ws.Range("L132").Value = 8080.799999999999
$ws.Range("M132").Value = -270.69812
$ws.Range("N132").Value = -13140.8

$ws.Range("H136").Value = 15874797
$ws.Range("I136").Value = 17922706
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 53768118
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -53765568
$ws.Range("N136").Value = -15600
